$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the date column (B2:B3) from a date-formatted number to text values
$ws.Range("B2:B3").NumberFormat = "@"
$ws.Range("B2").Value = "12/04/2022"
$ws.Range("B3").Value = "12/04/2022"

# Update the selected/active cell in the sheet view
$ws.Range("D6").Select() | Out-Null
